# Fill in the Plan (D) and Actual (E) values for the daily rows, and the
# starting Status carry-over (G5). Columns F and G contain formulas that
# will recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Starting carry-over value in G5
$ws.Range("G5").Value = -2378

# Row => (Plan D, Actual E)
$data = @{
    6  = @(110, 105)
    7  = @(110, 113)
    8  = @(110, 108)
    9  = @(110, 126)
    12 = @(144, 129)
    13 = @(110, 120)
    14 = @(110, 106)
    15 = @(110, 107)
    16 = @(110, 106)
    19 = @(144, 124)
    20 = @(110, 143)
    21 = @(110, 115)
    22 = @(110, 134)
    23 = @(110, 87)
    26 = @(110, 0)
    27 = @(110, 0)
    28 = @(110, 0)
    29 = @(110, 0)
    30 = @(110, 0)
    33 = @(110, 0)
    34 = @(41, 0)
}

foreach ($row in $data.Keys) {
    $plan = $data[$row][0]
    $actual = $data[$row][1]
    $ws.Cells.Item($row, 4).Value = $plan
    $ws.Cells.Item($row, 5).Value = $actual
}

$excel.CalculateFull()
